# Update countries & provincias Spain
# - Refresh case counts for Suiza, Costa Rica, Maldivas and Malaui.
# - Costa Rica / Maldivas overtake their neighbours (Crucero / Vietnam) in
#   total cases, so the table (sorted descending by "Casos totales") is
#   re-sorted to reflect the new ranking.
# - Bump the "last updated" timestamp shown in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$dataRange = $ws.Range("A4:H216")
$countryCol = $ws.Range("A4:A216")

# xlWhole = 1 -> match the whole cell content, not a partial substring
# (avoids e.g. "Costa Rica" accidentally matching "Costa de Marfil").
$xlWhole = 1

# --- Suiza: recovered / deaths-today / deaths updated ---
$row = $countryCol.Find("Suiza", [Type]::Missing, [Type]::Missing, $xlWhole).Row
$ws.Cells.Item($row, 5).Value2 = 5091
$ws.Cells.Item($row, 7).Value2 = 17
$ws.Cells.Item($row, 8).Value2 = 1716

# --- Costa Rica: total / new / active / recovered updated ---
$row = $countryCol.Find("Costa Rica", [Type]::Missing, [Type]::Missing, $xlWhole).Row
$ws.Cells.Item($row, 2).Value2 = 713
$ws.Cells.Item($row, 3).Value2 = 8
$ws.Cells.Item($row, 4).Value2 = 323
$ws.Cells.Item($row, 5).Value2 = 384

# --- Maldivas: total / new / recovered updated ---
$row = $countryCol.Find("Maldivas", [Type]::Missing, [Type]::Missing, $xlWhole).Row
$ws.Cells.Item($row, 2).Value2 = 278
$ws.Cells.Item($row, 3).Value2 = 28
$ws.Cells.Item($row, 5).Value2 = 261

# --- Malaui: active / recovered updated ---
$row = $countryCol.Find("Malaui", [Type]::Missing, [Type]::Missing, $xlWhole).Row
$ws.Cells.Item($row, 4).Value2 = 7
$ws.Cells.Item($row, 5).Value2 = 26

# Re-sort the whole table descending by "Casos totales" (col B) now that
# Costa Rica and Maldivas outrank the countries just above them.
$dataRange.Sort($ws.Range("B4:B216"), 2)

# --- Update the "datos actualizados" timestamp ---
$ws.Range("A1").Value2 = "Datos actualizados a 29 de Abril de 2020 a las 21:22"
